$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force cells to be treated as text so values are stored exactly
# (matching the original inlineStr representation) rather than being
# auto-converted to numbers/percentages by Excel.
$cells = @("D2","E2","D3","E3","D4","E4","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","E14","E15","D16","E16","D17","E17","D18","E18","E19","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D26","E26","E27","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","E47")
foreach ($cellRef in $cells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "280.96"
$ws.Range("E2").Value = "5.71%"
$ws.Range("D3").Value = "26.85"
$ws.Range("E3").Value = "0.03%"
$ws.Range("D4").Value = "4.942"
$ws.Range("E4").Value = "5.14%"
$ws.Range("E5").Value = "5.33%"
$ws.Range("D6").Value = "6.983"
$ws.Range("E6").Value = "3.65%"
$ws.Range("D7").Value = "3.353"
$ws.Range("E7").Value = "5.86%"
$ws.Range("D8").Value = "0.8859"
$ws.Range("E8").Value = "4.12%"
$ws.Range("D9").Value = "1.002"
$ws.Range("E9").Value = "10.35%"
$ws.Range("D10").Value = "0.1487"
$ws.Range("E10").Value = "5.71%"
$ws.Range("D11").Value = "0.05159"
$ws.Range("E11").Value = "1.68%"
$ws.Range("D12").Value = "0.07407"
$ws.Range("E12").Value = "4.46%"
$ws.Range("D13").Value = "0.03102"
$ws.Range("E13").Value = "-1.47%"
$ws.Range("E14").Value = "0.30%"
$ws.Range("E15").Value = "2.15%"
$ws.Range("D16").Value = "0.0006300"
$ws.Range("E16").Value = "3.81%"
$ws.Range("D17").Value = "0.006037"
$ws.Range("E17").Value = "-1.31%"
$ws.Range("D18").Value = "3.508"
$ws.Range("E18").Value = "1.51%"
$ws.Range("E19").Value = "5.67%"
$ws.Range("D21").Value = "0.1330"
$ws.Range("E21").Value = "3.79%"
$ws.Range("D22").Value = "3.940"
$ws.Range("E22").Value = "-3.70%"
$ws.Range("D23").Value = "0.04344"
$ws.Range("E23").Value = "2.38%"
$ws.Range("D24").Value = "0.001179"
$ws.Range("E24").Value = "0.04%"
$ws.Range("D25").Value = "0.003685"
$ws.Range("E25").Value = "-9.23%"
$ws.Range("D26").Value = "0.0001200"
$ws.Range("E26").Value = "-0.10%"
$ws.Range("E27").Value = "0.68%"
$ws.Range("D40").Value = "0.04098"
$ws.Range("E40").Value = "4.35%"
$ws.Range("D41").Value = "0.006646"
$ws.Range("E41").Value = "58.52%"
$ws.Range("D42").Value = "0.1178"
$ws.Range("E42").Value = "5.63%"
$ws.Range("D43").Value = "0.002360"
$ws.Range("E43").Value = "11.76%"
$ws.Range("D44").Value = "0.01306"
$ws.Range("E44").Value = "13.64%"
$ws.Range("D45").Value = "0.00005258"
$ws.Range("E45").Value = "2.53%"
$ws.Range("E47").Value = "814.71%"
